$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for columns C and D (td_sim_1 / record_atd)
$updates = @{
    3  = 187
    5  = 113
    7  = 68
    9  = 49
    11 = 28
    13 = 81
    15 = 48
    17 = 59
    19 = 117
    21 = 15
    23 = 242
    25 = 93
    27 = 289
    29 = 106
    31 = 373
    33 = 27
    35 = 34
    36 = 37
    38 = 79
    40 = 111
    42 = 160
    43 = 47
    45 = 39
    47 = 287
    49 = 57
    51 = 5
    53 = 1073
    55 = 32
    57 = 50
    59 = 134
    61 = 36
    63 = 44
    65 = 14
}

foreach ($row in $updates.Keys) {
    $newVal = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $newVal

    # Column D mirrors column C, except row 42 where D is the average of B and C
    if ($row -eq 42) {
        $ws.Cells.Item($row, 4).Value = 96.5
    } else {
        $ws.Cells.Item($row, 4).Value = $newVal
    }
}

# Row 66: C66 is the average of the updated td_sim_1 column (C2:C65)
$ws.Cells.Item(66, 3).Value = 125.2727272727273
